$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Add the new time-tracking rows (62-67) plus the new totals row (68).
#    We copy existing, already-formatted rows as templates so the new cells
#    pick up the same number formats / alignment / merges as the rest of the
#    table, then overwrite the template values with the real data.
# ---------------------------------------------------------------------------

# Row 62 looks like a "new day" row (date in col A, merged task in B:E,
# start/end time in F/G, minutes in I) -> use row 53 as the template.
$ws.Range("A53:G53").Copy($ws.Range("A62:G62"))
$ws.Range("I53").Copy($ws.Range("I62"))

# Row 63 is also a "new day" row -> same template.
$ws.Range("A53:G53").Copy($ws.Range("A63:G63"))
$ws.Range("I53").Copy($ws.Range("I63"))

# Rows 64-67 continue the same day (no date in col A) -> use row 54 as the
# template (merged task in B:E, start/end time in F/G, minutes in I, no A).
$ws.Range("B54:G54").Copy($ws.Range("B64:G64"))
$ws.Range("I54").Copy($ws.Range("I64"))

$ws.Range("B54:G54").Copy($ws.Range("B65:G65"))
$ws.Range("I54").Copy($ws.Range("I65"))

$ws.Range("B54:G54").Copy($ws.Range("B66:G66"))
$ws.Range("I54").Copy($ws.Range("I66"))

$ws.Range("B54:G54").Copy($ws.Range("B67:G67"))
$ws.Range("I54").Copy($ws.Range("I67"))

# Row 68 is the block-total row (minutes / h-min text / running grand total)
# -> use row 60 as the template.
$ws.Range("J60:L60").Copy($ws.Range("J68:L68"))

# --- Fill in the actual values -------------------------------------------

# Row 62: 2020-11-16, 14:00-16:00, Projekt-Handbuch Erstellung und Formatierung
$ws.Range("A62").Value = 44151
$ws.Range("B62").Value = "Projekt-Handbuch Erstellung und Formatierung"
$ws.Range("F62").Value = 0.58333333333333337
$ws.Range("G62").Value = 0.66666666666666663
$ws.Range("I62").Value = 120

# Row 63: 2020-11-21, 19:50-20:25, Projekt Handbuch aktualisiert
$ws.Range("A63").Value = 44156
$ws.Range("B63").Value = "Projekt Handbuch aktualisiert"
$ws.Range("F63").Value = 0.82638888888888884
$ws.Range("G63").Value = 0.85069444444444453
$ws.Range("I63").Value = 35

# Row 64: 20:25-20:30, Projektauftrags-Formular aktualisiert
$ws.Range("B64").Value = "Projektauftrags-Formular aktualisiert"
$ws.Range("F64").Value = 0.85069444444444453
$ws.Range("G64").Value = 0.85416666666666663
$ws.Range("I64").Value = 5

# Row 65: 20:30-20:35, Kurzbericht erstellt und dokumentiert (reuses existing text)
$ws.Range("B65").Value = "Kurzbericht erstellt und dokumentiert "
$ws.Range("F65").Value = 0.85416666666666663
$ws.Range("G65").Value = 0.85763888888888884
$ws.Range("I65").Value = 5

# Row 66: 20:35-20:40, Projekt auf Github raufgeladen
$ws.Range("B66").Value = "Projekt auf Github raufgeladen"
$ws.Range("F66").Value = 0.85763888888888884
$ws.Range("G66").Value = 0.86111111111111116
$ws.Range("I66").Value = 5

# Row 67: 20:40-20:45, Zeitmanagement laufend aktualisiert
$ws.Range("B67").Value = "Zeitmanagement laufend aktualisiert"
$ws.Range("F67").Value = 0.86111111111111116
$ws.Range("G67").Value = 0.86458333333333337
$ws.Range("I67").Value = 5

# Row 68: block totals
$ws.Range("J68").Value = 175
$ws.Range("K68").Value = "2h 55min"
$ws.Range("L68").Value = "33h"

# ---------------------------------------------------------------------------
# 2) Bold the "Summe (h)" column (header + grand-total figures).
# ---------------------------------------------------------------------------
$ws.Range("L1").Font.Bold = $true
$ws.Range("L50").Font.Bold = $true
$ws.Range("L60").Font.Bold = $true
$ws.Range("L68").Font.Bold = $true

# ---------------------------------------------------------------------------
# 3) Page setup (paper size / orientation) for printing.
# ---------------------------------------------------------------------------
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# ---------------------------------------------------------------------------
# 4) Update the selection so it reflects where the author ended up after
#    adding the new rows above.
# ---------------------------------------------------------------------------
$ws.Range("N74").Select()
